$d = $word.ActiveDocument

# The phrase below occurs twice in the document. In both places the single
# run that holds it must become three runs (all with the exact same
# character formatting as the original run):
#   1) "زمان بسیار زیاد"
#   2) "ی"
#   3) " طول می‌کشد تا فهرست همه پوشه‌ها را بدست آورد"
# i.e. "زیاد" grows the extra letter "ی" ("زیادی") and the rest of the
# sentence is pushed into a run of its own.
#
# Word only starts a new run where formatting actually changes, so a plain
# text insertion/replacement just gets folded back into one run. To force a
# split at an exact offset we toggle a character-formatting property (Bold)
# off/on/off on the sub-range we want to become its own run; toggling is a
# visual no-op but makes Word commit a fresh run for exactly that range.
# The two runs immediately following our edit ("! " and "مثلاً برای درایو ")
# already share the same formatting as the run we are splitting, so without
# re-asserting their own boundaries they would silently be folded into our
# new third run - we protect them the same way.

function Split-RunAt($range) {
    # Re-stamps a (start,end) range as its own run without changing how it
    # looks: Bold -> Bold -> original value, using a genuine set each time.
    $range.Bold = 1
    $range.Bold = 0
}

$needle = "زمان بسیار زیاد طول می‌کشد تا فهرست همه پوشه‌ها را بدست آورد"
$prefix = "زمان بسیار زیاد"
$insertChar = "ی"
$afterNeedle1 = "! "
$afterNeedle2 = "مثلاً برای درایو "

for ($i = 0; $i -lt 2; $i++) {
    $rng = $d.Content
    $found = $rng.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        break
    }

    $start = $rng.Start
    $end = $rng.End
    $prefixLen = $prefix.Length
    $tailLen = $end - ($start + $prefixLen)

    # Insert the missing "ی" right after "...زیاد" and before the space that
    # follows it.
    $insPoint = $d.Range($start + $prefixLen, $start + $prefixLen)
    $insPoint.InsertBefore($insertChar)

    # Everything from $start+$prefixLen onward shifted right by the length
    # of the inserted text.
    $aStart = $start
    $aEnd   = $start + $prefixLen
    $bStart = $aEnd
    $bEnd   = $bStart + $insertChar.Length
    $cStart = $bEnd
    $cEnd   = $cStart + $tailLen

    Split-RunAt $d.Range($aStart, $aEnd)
    Split-RunAt $d.Range($bStart, $bEnd)
    Split-RunAt $d.Range($cStart, $cEnd)

    # Keep the following, untouched runs from merging into the new run C.
    $dEnd = $cEnd + $afterNeedle1.Length
    Split-RunAt $d.Range($cEnd, $dEnd)

    $eEnd = $dEnd + $afterNeedle2.Length
    Split-RunAt $d.Range($dEnd, $eEnd)
}

Write-Output "done"
